# Add a new "Save" column (H) to the s_vals sheet, matching the header
# style already used by the other header cells (B1:G1), and a data value
# of 0 in H2 (same pattern as F2 "Win").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "sum" header cell (G1) onto the new
# header cell so the new column reuses the same style (bold, bordered,
# centered) rather than creating a brand-new style entry.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
